# cc_mech_calcs.xlsx - "updated all bolt-on components to have FS of 1.5"
#
# Reduce the plate thickness (t, cell H3) on the Injector Plate and Aft
# Closeout sheets so that their resulting factor of safety (FS) comes out
# to ~1.5, matching the rest of the bolt-on hardware. Also add two
# explanatory footnote rows on the Overview sheet, bold the Butt Weld
# FS cell to match the other FS cells in the workbook, and leave the
# various sheet selections where the author last left them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Injector Plate: t (H3) 0.5 -> 0.425  =>  FS (C16) ~= 1.4997 (~1.5)
# ---------------------------------------------------------------------
$wsInjector = $wb.Worksheets.Item("Injector Plate")
$wsInjector.Range("H3").Value = 0.425

# ---------------------------------------------------------------------
# Aft Closeout: t (H3) 0.5 -> 0.31  =>  FS (C22) ~= 1.5078 (~1.5)
# ---------------------------------------------------------------------
$wsAft = $wb.Worksheets.Item("Aft Closeout")
$wsAft.Range("H3").Value = 0.31

# ---------------------------------------------------------------------
# Overview: two new footnote rows under the Roark's reference note,
# defining the ta/rb moment abbreviations used elsewhere in the sheet.
# Written B14 then B13 so the shared-string table and references line
# up with the authored edit.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B14").Value = "ta = tangential moment at the outer edge"
$wsOverview.Range("B13").Value = "rb = radial moment at the inner edge"
$wsOverview.Range("B13:B14").Font.Size = 9

# ---------------------------------------------------------------------
# Butt Weld: bold the SF_t result cell to match the other FS cells
# (Flange C22/C35, Injector Plate C16, Aft Closeout C22, Proof Test
# Caps C16/C21) across the workbook.
# ---------------------------------------------------------------------
$wsButtWeld = $wb.Worksheets.Item("Butt Weld")
$wsButtWeld.Range("C11").Font.Bold = $true

# ---------------------------------------------------------------------
# Restore each sheet's last active selection. Activate each sheet to
# set its own selection, then finish back on Overview so it stays the
# tab that is shown when the workbook is reopened.
# ---------------------------------------------------------------------
$wsInjector.Activate()
$wsInjector.Range("H4").Select() | Out-Null

$wsAft.Activate()
$wsAft.Range("H4").Select() | Out-Null

$wsButtWeld.Activate()
$wsButtWeld.Range("C11").Select() | Out-Null

$wsOverview.Activate()
$wsOverview.Range("E17").Select() | Out-Null
